$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L2").Value = "[0.36740789940421664, 0.45943669909380086]"
$ws.Range("M2").Value = [double]"1.443289932012704e-14"
$ws.Range("N2").Value = [double]"1.443289932012704e-14"
$ws.Range("T2").Value = "[0.5407944418823556, 0.5898279023085289]"
